$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were removed in the diff
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E6").ClearContents()

# Update values that changed in the diff
$ws.Range("C7").Value = 0.9006569003772169
$ws.Range("E7").Value = 0.4141666450523163

$ws.Range("C8").Value = 1.298949644080372
$ws.Range("E8").Value = 0.7446484179501223

$ws.Range("C9").Value = 1.485127130420993
$ws.Range("E9").Value = 0.8988967199517361

$ws.Range("C10").Value = 1.938044824544427
$ws.Range("E10").Value = 1.197301207077017

$ws.Range("C11").Value = 1.730502563828185
$ws.Range("E11").Value = 1.20920901052266

$ws.Range("C12").Value = 2.211325510218898
$ws.Range("E12").Value = 1.513838358900466

$ws.Range("C13").Value = 1.095903126316466
$ws.Range("E13").Value = 1.063472944477306

$ws.Range("C14").Value = 0.8137456736830195
$ws.Range("E14").Value = 1.30966355756772

$ws.Range("C15").Value = -1.434438137829841
$ws.Range("E15").Value = 0.8159375071586261

$ws.Range("C16").Value = 1.85385197842538
$ws.Range("E16").Value = 1.2808239555127

$ws.Range("C17").Value = -0.6079479926716203
$ws.Range("E17").Value = 0.8021760422591839

$ws.Range("C18").Value = -0.06520462171909491
$ws.Range("E18").Value = 0.7367476213790747

$ws.Range("C19").Value = 0.5869668956646645
$ws.Range("E19").Value = 0.8208952814083625
